$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.297.27'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '2.062.32'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("E4").Value = '  +0.00%  '
$c = $ws.Range("D5")
$c.Value = "'233.25"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.12%  '
$c = $ws.Range("D6")
$c.Value = "'0.623"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.15%  '
$ws.Range("E7").Value = '  +0.04%  '
$c = $ws.Range("D8")
$c.Value = "'56.84"
$c.Style = "Normal"
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("E10").Value = '  +0.48%  '
$ws.Range("E11").Value = '  +0.87%  '
$ws.Range("D12").Value = '2.365.73'
$ws.Range("E12").Value = '  -0.18%  '
$c = $ws.Range("D13")
$c.Value = "'14.40"
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("E14").Value = '  -0.69%  '
$c = $ws.Range("D15")
$c.Value = "'0.772"
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.31%  '
$c = $ws.Range("D16")
$c.Value = "'5.14"
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.94%  '
$ws.Range("D17").Value = '2.061.75'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").Value = '37.240.38'
$ws.Range("E18").Value = '  -0.53%  '
$c = $ws.Range("D19")
$c.Value = "'6.37"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +3.78%  '
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("D21").Value = '0.0₃0811'
$ws.Range("E21").Value = '  +0.22%  '
$c = $ws.Range("D22")
$c.Value = "'226.17"
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("E23").Value = '  +0.02%  '
$c = $ws.Range("D24")
$c.Value = "'2.42"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.99%  '
$c = $ws.Range("D25")
$c.Value = "'2.40"
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.07%  '
$c = $ws.Range("D26")
$c.Value = "'166.18"
$c.Style = "Normal"
$ws.Range("E26").Value = '  +1.58%  '
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("E28").Value = '  +2.41%  '
$c = $ws.Range("D29")
$c.Value = "'19.03"
$c.Style = "Normal"
$ws.Range("E30").Value = '  -2.41%  '
$ws.Range("E31").Value = '  -0.73%  '
$c = $ws.Range("D32")
$c.Value = "'4.47"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.34%  '
$c = $ws.Range("D33")
$c.Value = "'4.60"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +3.57%  '
$c = $ws.Range("D34")
$c.Value = "'0.0616"
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.42%  '
$ws.Range("E35").Value = '  -4.81%  '
$ws.Range("E36").Value = '  -0.01%  '
$c = $ws.Range("D37")
$c.Value = "'1.79"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.50%  '
$c = $ws.Range("D38")
$c.Value = "'3.20"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -4.32%  '
$c = $ws.Range("D39")
$c.Value = "'5.52"
$c.Style = "Normal"
$ws.Range("E39").Value = '  -5.25%  '
$ws.Range("E40").Value = '  -0.80%  '
$ws.Range("D41").Value = '1.471.79'
$ws.Range("E41").Value = '  +0.48%  '
$c = $ws.Range("D42")
$c.Value = "'95.91"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D43")
$c.Value = "'0.0213"
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.98%  '
$ws.Range("E44").Value = '  +3.44%  '
$c = $ws.Range("D45")
$c.Value = "'4.30"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -3.95%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D46")
$c.Value = "'0.0928"
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.91%  '
$ws.Range("E47").Value = '  -0.35%  '
$c = $ws.Range("D48")
$c.Value = "'15.12"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -6.04%  '
$ws.Range("E49").Value = '  -1.49%  '
$c = $ws.Range("D50")
$c.Value = "'2.95"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.69%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$c = $ws.Range("D51")
$c.Value = "'44.23"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.13%  '
